# Added Hungary FC Test data
#
# Mirrors the "Slovakia" market tab (the last of the existing per-country
# sheets) into a brand-new "Hungary" tab positioned right after it, then
# swaps in the Hungary-specific strings for the two data cells that differ
# between markets (B2 = market name, B4 = NGC ticket reference).

$wb = $excel.ActiveWorkbook

# Slovakia is the template every market sheet is cloned from - same
# layout/styles, only B2 and B4 change between countries.
$slovakia = $wb.Worksheets.Item("Slovakia")

# Copy Slovakia and drop the clone immediately after it.
$slovakia.Copy($null, $slovakia)
$hungary = $wb.ActiveSheet
$hungary.Name = "Hungary"

# Market-specific values for the new sheet.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3590/T3615"

# Restore Slovakia to an unselected/default view state (it's no longer the
# active tab) and leave Hungary, the newly added sheet, active & selected.
$slovakia.Cells.Select()
$hungary.Activate()
$hungary.Range("B2:B4").Select()
